# Webprogrammierung_2.pptx - "Was möchten wir verbessern?" slide:
# Add a new bullet "Kontaktinformationen bei Gesuchen" right after the
# "Bessere Lösung zur Anlage eines Gesuches" bullet (and before the
# trailing empty paragraph), matching the commit
# "Neu Laden der Seite bei Anlage Gesuch".

$p = $ppt.ActivePresentation

# --- locate the bullet list shape that contains the "Bessere Lösung..." line ---
$targetShape = $null
$targetParaIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if (-not $sh.HasTextFrame) { continue }

        $tr = $sh.TextFrame.TextRange
        for ($pi = 1; $pi -le 50; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Start -eq 0) { break }

            if ($para.Text -like "*Anlage eines Gesuches*") {
                $targetShape = $sh
                $targetParaIndex = $pi
            }

            if ($para.Length -eq 0 -and $pi -gt 1 -and $para.Text -eq "") {
                # keep scanning a couple more in case of false positives, but a
                # zero-length paragraph normally marks the end of the body text
            }
        }
    }
}

if ($targetShape -eq $null) {
    Write-Host "Could not find the 'Bessere Loesung ... Gesuches' bullet."
} else {
    $tr = $targetShape.TextFrame.TextRange

    # Re-fetch the matched paragraph range and split a new paragraph right
    # after it (carriage return = new paragraph break, inherits the bullet's
    # pPr automatically, just like typing Enter at the end of the line).
    $bulletPara = $tr.Paragraphs($targetParaIndex, 1)
    $null = $bulletPara.InsertAfter([char]13 + "Kontaktinformationen ")

    # The freshly created paragraph is now the next one; append the second
    # run onto it so the line reads "Kontaktinformationen bei Gesuchen"
    # split across two runs, matching the authored edit.
    $newPara = $tr.Paragraphs($targetParaIndex + 1, 1)
    $null = $newPara.InsertAfter("bei Gesuchen")
}
